$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row permutation: new row $r gets the old D/J/K/L/M/P values that used to
# live in row $mapping[$r]. Captured from the authoritative diff (a pure
# reshuffle of the weekly price rows onto different dates/rows).
$mapping = @{
    2 = 15;  3 = 29;  4 = 8;   5 = 30;  6 = 35;  7 = 19;  8 = 21;  9 = 20;
    10 = 28; 11 = 16; 12 = 34; 13 = 2;  14 = 11; 15 = 3;  16 = 31; 17 = 33;
    18 = 39; 19 = 24; 20 = 4;  21 = 13; 22 = 36; 23 = 10; 24 = 37; 25 = 12;
    26 = 22; 27 = 14; 28 = 9;  29 = 26; 30 = 6;  31 = 23; 32 = 7;  33 = 25;
    34 = 18; 35 = 17; 36 = 27; 37 = 5;  38 = 32; 39 = 38
}

$cols = @("D", "J", "K", "L", "M", "P")

# Snapshot the original values for every touched column/row BEFORE any
# writes happen, since the mapping is a permutation (every row is both a
# source and a destination) and writes must not clobber a value that is
# still needed as someone else's source.
$original = @{}
foreach ($col in $cols) {
    for ($row = 2; $row -le 39; $row++) {
        $original["$col$row"] = $ws.Range("$col$row").Value2
    }
}

foreach ($row in $mapping.Keys) {
    $srcRow = $mapping[$row]
    foreach ($col in $cols) {
        $ws.Range("$col$row").Value2 = $original["$col$srcRow"]
    }
}
